# The "date" placeholder in the right-hand signature cell is preceded by a
# long run of spaces, all inside a single run:
#
#   <w:r><w:rPr>...</w:rPr>
#     <w:t xml:space="preserve">                               << Date >></w:t>
#   </w:r>
#
# The target edit splits that single run into three runs (same rPr on all
# three) and inserts 9 extra spaces right before "<< Date >>":
#
#   run 1: the original 31 spaces
#   run 2: 9 new spaces
#   run 3: "<< Date >>" (no leading space -> no xml:space="preserve")
#
# Locate the target occurrence of "<< Date >>" by searching for the one
# that is preceded by a run of nothing-but-whitespace (there is another,
# earlier "<< Date >>" in the table that is NOT preceded by spaces, so this
# distinguishes the two).

$d = $word.ActiveDocument

$precedingLen = 31
$targetStart = -1

$search = $d.Content
while ($search.Find.Execute("<< Date >>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $matchStart = $search.Start

    if ($matchStart -ge $precedingLen) {
        $preceding = $d.Range($matchStart - $precedingLen, $matchStart).Text
        if ($preceding.Trim().Length -eq 0) {
            $targetStart = $matchStart
        }
    }

    $search.Collapse(0)
}

if ($targetStart -ge 0) {
    # Insert 9 additional spaces immediately before "<< Date >>".
    $insertPoint = $d.Range($targetStart, $targetStart)
    $insertPoint.InsertBefore("         ")

    # The 9 freshly inserted spaces now occupy [$targetStart, $targetStart+9).
    # Toggling a direct character-formatting property on that sub-range and
    # then immediately back off forces Word to carve it out into its own
    # run (with rPr identical to its neighbours), splitting what was one
    # run into the three runs the target markup expects.
    $newSpaces = $d.Range($targetStart, $targetStart + 9)
    $newSpaces.Bold = 1
    $newSpaces.Bold = 0
}

Write-Host "targetStart=$targetStart"
